$d = $word.ActiveDocument

# 1. Replace the document body content with the full target structure
#    (adds the long narrative runs, empty paragraphs, page-break run,
#    and the _GoBack bookmark; also drops the stray paragraph-mark
#    language formatting that was on the original lone paragraph).
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>dsadf</w:t></w:r><w:r><w:t xml:space="preserve">алекое будущее. На основе ИИ была создана виртуальная реальность типа матрицы своего рода платформа создающая миры и контролирующая их функционал. Для того что бы сделать мир виртуальной реальности красивым, люди пошли на нарушение некоторых правил использования ИИ. Нужно было создать подобие реальности наполненной людьми способными испытывать чувства. Для этого человеческий разум копировался и создавались клоны людей в виртуальном мире считающие себя живыми людьми. Главный ИИ менял им воспоминания в соответствии с текущей реальностью. Личности для игр покупались фирмами как товар на который каждый человек имел авторские права. Кто-то приходил в виртуальные миры - отвлечься от реалий мира, просто жил там, но были и те кто вымещал на виртуальном мире свою злобу. Как правило это были очень богатые люди способные заплатить за убийство созданных клонов. Убитый клон стирался из системы навсегда. Главный герой сатоши был создателем проекта виртуальной реальности. Он был смертельно болен и перед смертью клонировал свою личность в виртуальный мир вписав в нее часть исходного кода платформы ИИ. Благодаря этому его нельзя было стереть. В виртуальном мире все продавалось за особую виртуальную валюту - биткоин: способности, прокачка уровня оружие и броня причем для поддержания иллюзии реальности у клонов было строго запрещено взаимодействие между клонами и людьми пришедшеми из реальности. Т.е. разговоры всех людей из реального мира изменялись системой. Есть два варианта: (как в </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>WOW</w:t></w:r><w:r><w:t xml:space="preserve"> две рассы не могу понимать друг друга) и вариант номер два - в случае попытки разговора о том что мир реальный - система блокирует способность человека из реального мира говорить. Так же были монстры созданные ИИ. Жители виртуального мира были слабы из-за нехватки финансов на прокачку, монстры же были сильными и представляли угрозу превращая их жизнь в ад. Допустим изначально игра позиционировалась как защита волшебного мира от этих монстров. Приходили люди из реального мира и сражались с монстрами защищая виртуальный. Люди из виртуального мира считали их великими волшебниками пришедшими из другого мира. Но появился злой волшебник подчиняющий себе монстров и стремящийся уничтожить мир. Допустим это жена Сатоши из реального мира которая после смерти мужа решила уничтожить мир который он создал(он проводил за работой больше времени чем с семьей)</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>главный герой живет в мире который зародился с помощью ИИ. Допустим изначально это было проектирование модоли матрицы. ИИ загружал в модель мира модели людей которым платились деньги за их копирование. Далее снималось шоу в котором данные люди попадали в определенные условия. Для контроля над системой существовали админы. Они использовали аватары для входа в систему. В созданном мире они были кем-то вроде инопланетян которые по легенде захватили мир и сделали из него полигон для разного рода игр на которые люди во внешнем мире делали ставки. Но первые программисты создавшие данный мир вписали в свои копии исходный код. Они обладали особой силой в данном мире. Так же в данном мире у людей были особые способности - скилы. Они не знали откуда у них эти силы. Все считали что они рождаются из портала. Так же когда люди гибли они рассыпались в пыль которая исчезала - и никто не знал почему. Все привыкли что это законы мира. Так же в этом мире были спецслужбы которые следили за тем что бы никто не знал правды - это порождения ИИ - типа агентов в матрице. Они охотились за теми кто бросал вызов системе. Людьми(клонами) восставшими против системы и пытающимися разрушить ядро и теми кто обладал исходным кодом. Но вычислить таких людей они могли только когда проявлялась их сила. ГГ сначала сражается в играх на выживание против монстров. Ему снятся странные сны про реальный мир о которых он никому не рассказывает. или монстры атакуют его город и он теряет близких после чего начинает войну с монстрами. Потом он находит команду таких же как он и вместе они отправляются в ядро узнать что же происходит с миром. На них так же начинают охоту админы считая их вирусом который может разрушить ядро. Постепенно раскрывается сюжет - гг вспоминает что он на самом деле лишь виртуальная копия. Так же они встречаются с дружелюбными программами самого ядра.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t>dsadf</w:t></w:r><w:r><w:t xml:space="preserve">алекое будущее. На основе ИИ была создана виртуальная реальность типа матрицы своего рода платформа создающая миры и контролирующая их функционал. Для того что бы сделать мир виртуальной реальности красивым, люди пошли на нарушение некоторых правил использования ИИ. Нужно было создать подобие реальности наполненной людьми способными испытывать чувства. Для этого человеческий разум копировался и создавались клоны людей в виртуальном мире считающие себя живыми людьми. Главный ИИ менял им воспоминания в соответствии с текущей реальностью. Личности для игр покупались фирмами как товар на который каждый человек имел авторские права. Кто-то приходил в виртуальные миры - отвлечься от реалий мира, просто жил там, но были и те кто вымещал на виртуальном мире свою злобу. Как правило это были очень богатые люди способные заплатить за убийство созданных клонов. Убитый клон стирался из системы навсегда. Главный герой сатоши был создателем проекта виртуальной реальности. Он был смертельно болен и перед смертью клонировал свою личность в виртуальный мир вписав в нее часть исходного кода платформы ИИ. Благодаря этому его нельзя было стереть. В виртуальном мире все продавалось за особую виртуальную валюту - биткоин: способности, прокачка уровня оружие и броня причем для поддержания иллюзии реальности у клонов было строго запрещено взаимодействие между клонами и людьми пришедшеми из реальности. Т.е. разговоры всех людей из реального мира изменялись системой. Есть два варианта: (как в </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>WOW</w:t></w:r><w:r><w:t xml:space="preserve"> две рассы не могу понимать друг друга) и вариант номер два - в случае попытки разговора о том что мир реальный - система блокирует способность человека из реального мира говорить. Так же были монстры созданные ИИ. Жители виртуального мира были слабы из-за нехватки финансов на прокачку, монстры же были сильными и представляли угрозу превращая их жизнь в ад. Допустим изначально игра позиционировалась как защита волшебного мира от этих монстров. Приходили люди из реального мира и сражались с монстрами защищая виртуальный. Люди из виртуального мира считали их великими волшебниками пришедшими из другого мира. Но появился злой волшебник подчиняющий себе монстров и стремящийся уничтожить мир. Допустим это жена Сатоши из реального мира которая после смерти мужа решила уничтожить мир который он создал(он проводил за работой больше времени чем с семьей)</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>главный герой живет в мире который зародился с помощью ИИ. Допустим изначально это было проектирование модоли матрицы. ИИ загружал в модель мира модели людей которым платились деньги за их копирование. Далее снималось шоу в котором данные люди попадали в определенные условия. Для контроля над системой существовали админы. Они использовали аватары для входа в систему. В созданном мире они были кем-то вроде инопланетян которые по легенде захватили мир и сделали из него полигон для разного рода игр на которые люди во внешнем мире делали ставки. Но первые программисты создавшие данный мир вписали в свои копии исходный код. Они обладали особой силой в данном мире. Так же в данном мире у людей были особые способности - скилы. Они не знали откуда у них эти силы. Все считали что они рождаются из портала. Так же когда люди гибли они рассыпались в пыль которая исчезала - и никто не знал почему. Все привыкли что это законы мира. Так же в этом мире были спецслужбы которые следили за тем что бы никто не знал правды - это порождения ИИ - типа агентов в матрице. Они охотились за теми кто бросал вызов системе. Людьми(клонами) восставшими против системы и пытающимися разрушить ядро и теми кто обладал исходным кодом. Но вычислить таких людей они могли только когда проявлялась их сила. ГГ сначала сражается в играх на выживание против монстров. Ему снятся странные сны про реальный мир о которых он никому не рассказывает. или монстры атакуют его город и он теряет близких после чего начинает войну с монстрами. Потом он находит команду таких же как он и вместе они отправляются в ядро узнать что же происходит с миром. На них так же начинают охоту админы считая их вирусом который может разрушить ядро. Постепенно раскрывается сюжет - гг вспоминает что он на самом деле лишь виртуальная копия. Так же они встречаются с дружелюбными программами самого ядра.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:sectPr><w:pgSz w:w="11906" w:h="16838"/><w:pgMar w:top="1134" w:right="850" w:bottom="1134" w:left="1701" w:header="708" w:footer="708" w:gutter="0"/><w:cols w:space="708"/><w:docGrid w:linePitch="360"/></w:sectPr></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Content.InsertXML($xml)

# 2. Normal style now carries explicit line spacing (256/auto), matching
#    the updated style definition in styles.xml.
$normal = $d.Styles("Normal")
$normal.ParagraphFormat.LineSpacingRule = 5
$normal.ParagraphFormat.LineSpacing = 12.8
